$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper: set the TextToDisplay of the hyperlink anchored at a
# given A1 address on a given worksheet (address stays untouched,
# only the cached display text / cell text changes).
# ---------------------------------------------------------------
function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
            break
        }
    }
}

# ---------------------------------------------------------------
# Helper: remove the hyperlink anchored at a given A1 address
# ---------------------------------------------------------------
function Remove-HyperlinkAt($ws, [string]$addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            break
        }
    }
}

# =================================================================
# Sheet "Overview": new handoff info for the two source files
# =================================================================
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A2").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
$wsOv.Range("B2").Value = "Ready for handoff"
$wsOv.Range("C2").Value = "Ready for handoff"
$wsOv.Range("D2").Value = "2016-50-17 22:50:13"

$wsOv.Range("A3").Value = "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"
$wsOv.Range("B3").Value = "Ready for handoff"
$wsOv.Range("C3").Value = "Ready for handoff"
$wsOv.Range("D3").Value = "2016-50-17 22:50:13"

Set-HyperlinkDisplay $wsOv '$A$2' "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
Set-HyperlinkDisplay $wsOv '$A$3' "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"

# =================================================================
# Sheet "zh-cn": new handoff info + drop Target/Handback columns
# =================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-17 22:50:10"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value = "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-17 22:50:10"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"

Set-HyperlinkDisplay $wsZh '$A$2' "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
Set-HyperlinkDisplay $wsZh '$D$2' "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$A$3' "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"
Set-HyperlinkDisplay $wsZh '$D$3' "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.zh-cn.xlf"

Remove-HyperlinkAt $wsZh '$F$2'
Remove-HyperlinkAt $wsZh '$G$2'
Remove-HyperlinkAt $wsZh '$F$3'
Remove-HyperlinkAt $wsZh '$G$3'
$wsZh.Range("F2:G3").Clear()

# =================================================================
# Sheet "de-de": new handoff info + drop Target/Handback columns
# =================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-17 22:50:13"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value = "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-17 22:50:13"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"

Set-HyperlinkDisplay $wsDe '$A$2' "681e51b2-76c5-4dee-8010-856d0a93b3e4.md"
Set-HyperlinkDisplay $wsDe '$D$2' "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$A$3' "ffffd2e03a40-0903-46d7-8121-c6594f46e066.md"
Set-HyperlinkDisplay $wsDe '$D$3' "681e51b2-76c5-4dee-8010-856d0a93b3e4.2b05a74fc1a434de9b3f96543d4c4e6d89313243.de-de.xlf"

Remove-HyperlinkAt $wsDe '$F$2'
Remove-HyperlinkAt $wsDe '$G$2'
Remove-HyperlinkAt $wsDe '$F$3'
Remove-HyperlinkAt $wsDe '$G$3'
$wsDe.Range("F2:G3").Clear()
